$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 373 (shifts existing rows 373-452 down to 375-454),
# shifting cells down only within the used column range so the sheet dimension/extent
# doesn't balloon out to the full row width.
$ws.Range("A373:T374").Insert(-4121)  # xlShiftDown

function Set-Row {
    param($r, $variedad, $calidad, $volumen, $pmin, $pmax, $pprom, $pkg)
    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 4).Value = 45015
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102005
    $ws.Cells.Item($r, 10).Value = "Naranja"
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $pmin
    $ws.Cells.Item($r, 15).Value = $pmax
    $ws.Cells.Item($r, 16).Value = $pprom
    $ws.Cells.Item($r, 17).Value = "`$/caja 15 kilos empedrada"
    $ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
    $ws.Cells.Item($r, 19).Value = $pkg
    $ws.Cells.Item($r, 20).Value = 15
}

Set-Row 373 "Valencia" "Primera" 100 13000 14000 13500 900
Set-Row 374 "Valencia" "Segunda" 50 11000 11000 11000 733
